$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.349533081054688
$ws.Range("B1").Value = 2.133534908294678
$ws.Range("C1").Value = 2.034304141998291
$ws.Range("D1").Value = 2.718106269836426
$ws.Range("E1").Value = 4.220511913299561
